# Update the Java stack-trace text inside the 'setConserveRatio' test-output run
# to match the line numbers/call chain produced after moving from POI 3.16 to 3.17
# (Fixed #253). Each block below is a small, uniquely-anchored Find & Replace so that
# unrelated text (e.g. the straight apostrophe in "Can't") is left untouched.

$d = $word.ActiveDocument

# 0: MImageImpl.java:148->162 / ImageServices.java:156->188
$old0 = '	at org.obeonetwork.m2doc.element.impl.MImageImpl.setConserveRatio(MImageImpl.java:148)
	at org.obeonetwork.m2doc.services.ImageServices.setConserveRatio(ImageServices.java:156)
'
$new0 = '	at org.obeonetwork.m2doc.element.impl.MImageImpl.setConserveRatio(MImageImpl.java:162)
	at org.obeonetwork.m2doc.services.ImageServices.setConserveRatio(ImageServices.java:188)
'
$found0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)
if (-not $found0) { throw "Replacement 0 failed to find its target text" }
Write-Host "Replace 0: $found0"

# 1: M2DocEvaluator.caseQuery M2DocEvaluator.java:476->540
$old1 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:476)
'
$new1 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:540)
'
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) { throw "Replacement 1 failed to find its target text" }
Write-Host "Replace 1: $found1"

# 2: TemplateSwitch.java:172->186
$old2 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:172)
'
$new2 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)
'
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) { throw "Replacement 2 failed to find its target text" }
Write-Host "Replace 2: $found2"

# 3: M2DocEvaluator.java:836->1038 / caseBlock M2DocEvaluator.java:1034->1254
$old3 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1034)
'
$new3 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)
'
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found3) { throw "Replacement 3 failed to find its target text" }
Write-Host "Replace 3: $found3"

# 4: TemplateSwitch.java:183->199
$old4 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)
'
$new4 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)
'
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $found4) { throw "Replacement 4 failed to find its target text" }
Write-Host "Replace 4: $found4"

# 5: caseTemplate block replaced by caseDocumentTemplate block (lines 31-35 hunk)
$old5 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:297)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)
'
$new5 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)
'
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $found5) { throw "Replacement 5 failed to find its target text" }
Write-Host "Replace 5: $found5"

# 6: tail of stack (doSwitch/generate/M2DocUtils/AbstractTemplatesTestSuite) renumbered
$old6 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:259)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:252)
	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:691)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:396)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:318)
'
$new6 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)
	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)
'
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
if (-not $found6) { throw "Replacement 6 failed to find its target text" }
Write-Host "Replace 6: $found6"

# 7: insert RunAfters.evaluate line before ParentRunner.runLeaf
$old7 = '	at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)
	at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)
'
$new7 = '	at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)
	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
	at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)
'
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
if (-not $found7) { throw "Replacement 7 failed to find its target text" }
Write-Host "Replace 7: $found7"

# 8: remove RunBefores.evaluate line
$old8 = '	at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)
'
$new8 = ''
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
if (-not $found8) { throw "Replacement 8 failed to find its target text" }
Write-Host "Replace 8: $found8"

# 9: insert extra ParentRunner/Suite frames after RunAfters.evaluate/ParentRunner.run
$old9 = '	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
'
$new9 = '	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.junit.runners.Suite.runChild(Suite.java:128)
	at org.junit.runners.Suite.runChild(Suite.java:27)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.junit.runners.Suite.runChild(Suite.java:128)
	at org.junit.runners.Suite.runChild(Suite.java:27)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
'
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
if (-not $found9) { throw "Replacement 9 failed to find its target text" }
Write-Host "Replace 9: $found9"

# 10: RemoteTestRunner line numbers renumbered
$old10 = '	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:459)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:675)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:382)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:192)
'
$new10 = '	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)
'
$found10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2)
if (-not $found10) { throw "Replacement 10 failed to find its target text" }
Write-Host "Replace 10: $found10"

